# Bugfix: complete data by using values both from codelist and data
# A row describing the "age" codelist's catch-all/placeholder code (99) was
# missing from the "Codelists" sheet - insert it at row 103 (pushing the
# existing rows down), then re-select that sheet/cell as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codelists")

# Insert a new row above the current row 103, shifting rows 103:111 down to 104:112
$ws.Rows.Item(103).Insert()

# Fill in the new row with the missing codelist entry
$ws.Range("A103").Value = "age"
$ws.Range("B103").Value = 99
$ws.Range("C103").Value = 99
$ws.Range("D103").Value = 99
$ws.Range("E103").Value = 99
$ws.Range("F103").Value = 99

# Make "Codelists" the active sheet/tab, matching the author's final selection
$ws.Activate()
$ws.Range("G103").Select()
